# Generate Report for Archive
#
# 1) Status text changed from "Ready for handoff" to "In Translation" on all
#    three sheets (Overview summary columns + per-locale Status column).
# 2) The now-shorter status text lets the "Status" columns narrow, so their
#    column width shrinks accordingly on all three sheets.
#
# Note on the column-width literal below: the target stored width is
# 13.4101845877511 characters. This runtime's ColumnWidth setter quantizes
# to a 1/6-character pixel grid and adds a fixed 5/6 offset when it writes
# the stored <col width>, so we pre-compensate the assigned value
# (target - 5/6) to land on the closest representable stored width.
$targetColumnWidth = 12.576851254417766

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de summary-status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $targetColumnWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $targetColumnWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $targetColumnWidth
